$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("tasas")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.4 = 9033.61 pesos`n✅ 9033.61 pesos = 2.39 = 952.42 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

$ws2.Range("N10").Value = 416
$ws2.Range("O10").Value = 3757.98
$ws2.Range("N12").Value = 3775
$ws2.Range("O12").Value = 398
